$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was "M") becomes "B" row with new metrics
$ws.Range("A2").Value = "B"
$ws.Range("B2").Value = 0.9205298013245033
$ws.Range("C2").Value = 0.972027972027972
$ws.Range("D2").Value = 0.9455782312925171
$ws.Range("E2").Value = 143

# Row 3 (was "B") becomes "M" row with new metrics
$ws.Range("A3").Value = "M"
$ws.Range("B3").Value = 0.948051948051948
$ws.Range("C3").Value = 0.8588235294117647
$ws.Range("D3").Value = 0.9012345679012346
$ws.Range("E3").Value = 85

# Row 4 (accuracy) - all four columns get same new accuracy value
$ws.Range("B4").Value = 0.9298245614035088
$ws.Range("C4").Value = 0.9298245614035088
$ws.Range("D4").Value = 0.9298245614035088
$ws.Range("E4").Value = 0.9298245614035088

# Row 5 (macro avg)
$ws.Range("B5").Value = 0.9342908746882257
$ws.Range("C5").Value = 0.9154257507198683
$ws.Range("D5").Value = 0.9234063995968758

# Row 6 (weighted avg)
$ws.Range("B6").Value = 0.9307902507623665
$ws.Range("C6").Value = 0.9298245614035088
$ws.Range("D6").Value = 0.9290466023966442
